# HOBO Data first analysis
#
# Applies the edits described by the commit:
#   1. UPDATED 03/04/2024 -> UPDATED 04/01/2024, and move the "_GoBack"
#      bookmark from the end of "BLUETOOTH ALWAYS ON = Off" into the
#      middle of the updated date (right after the new day digits).
#   2. Merge the "Connect to data logger " / "by pushing the button..."
#      runs into a single run (no visible text change).
#   3. Merge the "**FIGURE O" / "UT BEST WAY..." runs into a single run
#      (no visible text change).
#   4. Rewrite the "Back up data..." bullet to reference GitHub and the
#      new HOBO naming convention.

$d = $word.ActiveDocument

function Split-RunAt($pos) {
    # Forces a run boundary at a collapsed position by briefly adding
    # and then removing a bookmark there - the split survives the
    # bookmark's removal.
    $tmpName = "zzTempSplitMarker"
    if ($d.Bookmarks.Exists($tmpName)) {
        $d.Bookmarks.Item($tmpName).Delete()
    }
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($tmpName, $r)
    $d.Bookmarks.Item($tmpName).Delete()
}

# ---------------------------------------------------------------------
# 1) Date line
# ---------------------------------------------------------------------

# Drop the _GoBack bookmark from its old location (end of the
# "BLUETOOTH ALWAYS ON = Off" bullet) - it is being relocated.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$datePara = $d.Paragraphs.Item(1).Range
$dateStart = $datePara.Start

# Sanity check on the text we expect to edit.
# "UPDATED 03/04/2024"
#  0123456789012345678
# Place the relocated bookmark between the new day ("01") and the
# slash before the year, i.e. right after "UPDATED 04/01".
$bmPos = $dateStart + ("UPDATED 04/01").Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Swap the month/day digits: 03/04 -> 04/01 (year is unchanged).
$d.Range($dateStart + 9, $dateStart + 10).Text = "4"
$d.Range($dateStart + 12, $dateStart + 13).Text = "1"

# Recreate the same run boundaries the live edit produced: "UPDATED 0"
# | "4" | "/0" | "1" | <bookmark> | "/2024".
Split-RunAt ($dateStart + 9)
Split-RunAt ($dateStart + 10)
Split-RunAt ($dateStart + 12)

# ---------------------------------------------------------------------
# 2) Merge "Connect to data logger " and "by pushing the button..."
#    into a single contiguous run.
# ---------------------------------------------------------------------
$connectText = "Connect to data logger by pushing the button on the logger while enabling "
$find = $d.Content.Find
$find.Execute($connectText, $true, $false, $false, $false, $false, $true, 1, $false)
$connectRange = $find.Parent
$connectRange.Text = ""
$connectRange.Collapse(1)
$connectRange.InsertBefore($connectText)

# ---------------------------------------------------------------------
# 3) Merge "**FIGURE O" and "UT BEST WAY..." into a single run.
# ---------------------------------------------------------------------
$figureText = "**FIGURE OUT BEST WAY TO SAVE DATA IN THE FIELD - excel sheet on the phone"
$find2 = $d.Content.Find
$find2.Execute($figureText, $true, $false, $false, $false, $false, $true, 1, $false)
$figureRange = $find2.Parent
$figureRange.Text = ""
$figureRange.Collapse(1)
$figureRange.InsertBefore($figureText)

# ---------------------------------------------------------------------
# 4) Rewrite the "Back up data..." bullet.
# ---------------------------------------------------------------------

# "Back up data into power bank (...). Label Data file as: " ->
# "Back up data into GitHub. Label Data file as: "
$oldLead = "Back up data into power bank (depending on how we figure out how to get it off phone. Otherwise will wait till service and email it). Label Data file as: "
$d.Content.Find.Execute($oldLead, $true, $false, $false, $false, $false, $true, 1, $false, "Back up data into GitHub. Label Data file as: ", 2)

# "LOGGERNO.SITENO.DATE.TIME" (still wrapped in the original
# gramStart/gramEnd proofErr pair) -> "HOBOmm.yyyy.LOGGERNO"
$d.Content.Find.Execute("LOGGERNO.SITENO.DATE.TIME", $true, $false, $false, $false, $false, $true, 1, $false, "HOBOmm.yyyy.LOGGERNO", 2)

# Trailing "." run (just after the old gramEnd) -> "_SITENO"
$d.Content.Find.Execute("HOBOmm.yyyy.LOGGERNO.", $true, $false, $false, $false, $false, $true, 1, $false, "HOBOmm.yyyy.LOGGERNO_SITENO", 2)
